$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$ws.Range("A44").Value = "Height"
$ws.Range("B44").Value = "float"
$ws.Range("C44").Value = $false
$ws.Range("D44").Value = $false
$ws.Range("E44").Value = $false
$ws.Range("F44").Value = $true
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = "Friend"
$ws.Range("I44").NumberFormat = "@"
$ws.Range("J44").Value = "模型高度"

[void]$ws.Range("J44").Select()
